$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(" Dubai (DSC)", " October 14 2020", "Capitals won by 13 runs", "Rajasthan Royals", "Delhi Capitals", "Riyan Parag$([char]0xA0)", "1", "2", "0", "0", "50.00"),
    @(" Dubai (DSC)", " October 22 2020", "Sunrisers won by 8 wickets (with 11 balls remaining)", "Rajasthan Royals", "Sunrisers Hyderabad", "Riyan Parag$([char]0xA0)", "20", "12", "2", "1", "166.66"),
    @(" Sharjah", " September 22 2020", "Royals won by 16 runs", "Rajasthan Royals", "Chennai Super Kings", "Riyan Parag$([char]0xA0)", "6", "4", "1", "0", "150.00"),
    @(" Dubai (DSC)", " September 30 2020", "KKR won by 37 runs", "Rajasthan Royals", "Kolkata Knight Riders", "Riyan Parag$([char]0xA0)", "1", "6", "0", "0", "16.66"),
    @(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Rajasthan Royals", "Kolkata Knight Riders", "Riyan Parag$([char]0xA0)", "0", "7", "0", "0", "0.00"),
    @(" Sharjah", " September 27 2020", "Royals won by 4 wickets (with 3 balls remaining)", "Rajasthan Royals", "Kings XI Punjab", "Riyan Parag$([char]0xA0)", "0", "2", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 11 2020", "Royals won by 5 wickets (with 1 ball remaining)", "Rajasthan Royals", "Sunrisers Hyderabad", "Riyan Parag$([char]0xA0)", "42", "26", "2", "2", "161.53"),
    @(" Abu Dhabi", " October 03 2020", "RCB won by 8 wickets (with 5 balls remaining)", "Rajasthan Royals", "Royal Challengers Bangalore", "Riyan Parag$([char]0xA0)", "16", "18", "1", "0", "88.88")
)

$startRow = 10
$endRow = $startRow + $data.Count - 1
$ws.Range("A$startRow`:K$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
